$d = $word.ActiveDocument

$replacements = @(
    @{old="288×8=2304"; new="205×4=820"},
    @{old="449×9=4041"; new="590×2=1180"},
    @{old="675×5=3375"; new="835×9=7515"},
    @{old="433×6=2598"; new="116×4=464"},
    @{old="541×4=2164"; new="673×7=4711"},
    @{old="250×8=2000"; new="910×8=7280"},
    @{old="883×2=1766"; new="509×4=2036"},
    @{old="211×2=422"; new="663×3=1989"},
    @{old="746×6=4476"; new="826×9=7434"},
    @{old="306×2=612"; new="156×7=1092"},
    @{old="385×6=2310"; new="345×9=3105"},
    @{old="163×8=1304"; new="735×5=3675"},
    @{old="755×8=6040"; new="657×8=5256"},
    @{old="538×7=3766"; new="586×2=1172"},
    @{old="255×7=1785"; new="771×2=1542"},
    @{old="625×8=5000"; new="705×6=4230"},
    @{old="587×9=5283"; new="694×7=4858"},
    @{old="817×6=4902"; new="430×3=1290"},
    @{old="292×7=2044"; new="951×3=2853"},
    @{old="679×5=3395"; new="511×3=1533"},
    @{old="982×4=3928"; new="749×8=5992"},
    @{old="692×8=5536"; new="364×9=3276"},
    @{old="881×3=2643"; new="283×4=1132"},
    @{old="638×9=5742"; new="125×4=500"},
    @{old="504×2=1008"; new="958×6=5748"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
